$wb = $excel.ActiveWorkbook

# ---- Sheet index 1 ----
$ws = $wb.Worksheets.Item(1)
# Row 17
$ws.Range("H17").Value = 1464.1666
$ws.Range("I17").Value = 1063
$ws.Range("J17").Value = 1500.6364
$ws.Range("K17").Value = 3189
$ws.Range("L17").Value = 4501.9092
$ws.Range("M17").Value = -3021
$ws.Range("N17").Value = -4837.9092
# Row 93
$ws.Range("H93").Value = 52271.5
$ws.Range("J93").Value = 52271.5
$ws.Range("L93").Value = 52271.5
$ws.Range("N93").Value = -57263.5
# Row 116
$ws.Range("H116").Value = 17863056
$ws.Range("I116").Value = 35718840
$ws.Range("J116").Value = 7271.2856
$ws.Range("K116").Value = 35718840
$ws.Range("L116").Value = 7271.2856
$ws.Range("M116").Value = -35715398
$ws.Range("N116").Value = -14155.2856
# Row 128
$ws.Range("H128").Value = 88917.60000000001
$ws.Range("J128").Value = 88917.60000000001
$ws.Range("L128").Value = 88917.60000000001
$ws.Range("N128").Value = -98877.60000000001
# Row 132
$ws.Range("H132").Value = 917.5714
$ws.Range("I132").Value = 958.15
$ws.Range("K132").Value = 2874.45
$ws.Range("M132").Value = -344.4499999999998
# Row 135
$ws.Range("H135").Value = 2500373.8
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
# Row 137
$ws.Range("H137").Value = 6492.9287
$ws.Range("I137").Value = 5834.1665
$ws.Range("J137").Value = 6987
$ws.Range("K137").Value = 17502.4995
$ws.Range("L137").Value = 20961
$ws.Range("M137").Value = -14952.4995
$ws.Range("N137").Value = -26061
# Row 141
$ws.Range("H141").Value = 6321.625
$ws.Range("I141").Value = 5178.8335
$ws.Range("K141").Value = 15536.5005
$ws.Range("M141").Value = -10356.5005

# ---- Sheet index 2 ----
$ws = $wb.Worksheets.Item(2)
# Row 2
$ws.Range("H2").Value = 50001772
$ws.Range("I2").Value = 1364.5714
$ws.Range("K2").Value = 1364.5714
$ws.Range("M2").Value = -1251.5714
# Row 32
$ws.Range("H32").Value = 4262605
$ws.Range("I32").Value = 4766963
$ws.Range("K32").Value = 4766963
$ws.Range("M32").Value = -4766676
# Row 47
$ws.Range("H47").Value = 21997
$ws.Range("J47").Value = 21997
$ws.Range("L47").Value = 21997
$ws.Range("N47").Value = -23447
# Row 61
$ws.Range("H61").Value = 27032138
$ws.Range("I61").Value = 2110.3809
$ws.Range("K61").Value = 2110.3809
$ws.Range("M61").Value = -1898.3809
# Row 74
$ws.Range("H74").Value = 18597
$ws.Range("I74").Value = 21987.291
$ws.Range("K74").Value = 21987.291
$ws.Range("M74").Value = -21113.291
# Row 77
$ws.Range("H77").Value = 18597
$ws.Range("I77").Value = 21987.291
$ws.Range("K77").Value = 109936.455
$ws.Range("M77").Value = -105568.455
# Row 116
$ws.Range("H116").Value = 50001772
$ws.Range("I116").Value = 1364.5714
$ws.Range("K116").Value = 1364.5714
$ws.Range("M116").Value = 929.4286
# Row 122
$ws.Range("H122").Value = 6959.6313
$ws.Range("I122").Value = 5911.875
$ws.Range("K122").Value = 17735.625
$ws.Range("M122").Value = -15285.625
# Row 136
$ws.Range("H136").Value = 27032138
$ws.Range("I136").Value = 2110.3809
$ws.Range("K136").Value = 6331.1427
$ws.Range("M136").Value = -3781.1427

# ---- Sheet index 3 ----
$ws = $wb.Worksheets.Item(3)
# Row 3
$ws.Range("H3").Value = 50001772
$ws.Range("I3").Value = 1364.5714
$ws.Range("K3").Value = 1364.5714
$ws.Range("M3").Value = -1250.5714
# Row 20
$ws.Range("H20").Value = 10421442
$ws.Range("I20").Value = 15155334
$ws.Range("J20").Value = 6879.6
$ws.Range("K20").Value = 15155334
$ws.Range("L20").Value = 6879.6
$ws.Range("M20").Value = -15155087
$ws.Range("N20").Value = -7373.6
# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
# Row 134
$ws.Range("H134").Value = 6763260.5
$ws.Range("I134").Value = 20835400
$ws.Range("K134").Value = 62506200
$ws.Range("M134").Value = -62503665

# ---- Sheet index 4 ----
$ws = $wb.Worksheets.Item(4)
# Row 7
$ws.Range("H7").Value = 93.72727
$ws.Range("I7").Value = 97.625
$ws.Range("K7").Value = 97.625
$ws.Range("M7").Value = 15.375
# Row 31
$ws.Range("H31").Value = 8095.3823
$ws.Range("I31").Value = 2252.3333
$ws.Range("K31").Value = 2252.3333
$ws.Range("M31").Value = -1957.3333
# Row 34
$ws.Range("H34").Value = 8095.3823
$ws.Range("I34").Value = 2252.3333
$ws.Range("K34").Value = 2252.3333
$ws.Range("M34").Value = -2050.3333
# Row 51
$ws.Range("H51").Value = 43612.5
$ws.Range("J51").Value = 43612.5
$ws.Range("L51").Value = 43612.5
$ws.Range("N51").Value = -45084.5
# Row 58
$ws.Range("H58").Value = 6948.3125
$ws.Range("I58").Value = 1792
$ws.Range("J58").Value = 10042.1
$ws.Range("K58").Value = 1792
$ws.Range("L58").Value = 10042.1
$ws.Range("M58").Value = -1589
$ws.Range("N58").Value = -10448.1
# Row 61
$ws.Range("H61").Value = 43612.5
$ws.Range("J61").Value = 43612.5
$ws.Range("L61").Value = 43612.5
$ws.Range("N61").Value = -44308.5
# Row 107
$ws.Range("H107").Value = 1858.7407
$ws.Range("I107").Value = 549.7143
$ws.Range("J107").Value = 2316.9
$ws.Range("K107").Value = 549.7143
$ws.Range("L107").Value = 2316.9
$ws.Range("M107").Value = 1370.2857
$ws.Range("N107").Value = -6156.9
# Row 136
$ws.Range("H136").Value = 6948.3125
$ws.Range("I136").Value = 1792
$ws.Range("J136").Value = 10042.1
$ws.Range("K136").Value = 5376
$ws.Range("L136").Value = 30126.3
$ws.Range("M136").Value = -2826
$ws.Range("N136").Value = -35226.3

# ---- Sheet index 5 ----
$ws = $wb.Worksheets.Item(5)
# Row 6
$ws.Range("H6").Value = 20833446
$ws.Range("I6").Value = 31250092
$ws.Range("K6").Value = 93750276
$ws.Range("M6").Value = -93750163
# Row 12
$ws.Range("H12").Value = 306.3125
$ws.Range("J12").Value = 71.57143000000001
$ws.Range("L12").Value = 214.71429
$ws.Range("N12").Value = -560.71429
# Row 32
$ws.Range("H32").Value = 166666750
$ws.Range("I32").Value = 500000060
$ws.Range("K32").Value = 1500000180
$ws.Range("M32").Value = -1499999897
# Row 131
$ws.Range("H131").Value = 2475.814
$ws.Range("I131").Value = 1819.0834
$ws.Range("J131").Value = 2730.0322
$ws.Range("K131").Value = 5457.2502
$ws.Range("L131").Value = 8190.096600000001
$ws.Range("M131").Value = -417.2502000000004
$ws.Range("N131").Value = -18270.0966
# Row 140
$ws.Range("H140").Value = 120192.82
$ws.Range("I140").Value = 155329.39
$ws.Range("K140").Value = 465988.17
$ws.Range("M140").Value = -460808.17

# ---- Sheet index 6 ----
$ws = $wb.Worksheets.Item(6)
# Row 97
$ws.Range("H97").Value = 645.65
$ws.Range("I97").Value = 573.3461
$ws.Range("K97").Value = 573.3461
$ws.Range("M97").Value = -77.34609999999998
# Row 102
$ws.Range("H102").Value = 2666.0278
$ws.Range("I102").Value = 2600.037
$ws.Range("J102").Value = 2864
$ws.Range("K102").Value = 2600.037
$ws.Range("L102").Value = 2864
$ws.Range("M102").Value = -978.0369999999998
$ws.Range("N102").Value = -6108
# Row 107
$ws.Range("H107").Value = 572131.9
$ws.Range("J107").Value = 297
$ws.Range("L107").Value = 297
$ws.Range("N107").Value = -4137
# Row 126
$ws.Range("H126").Value = 26321364
$ws.Range("I126").Value = 166668200
$ws.Range("J126").Value = 6333.125
$ws.Range("K126").Value = 500004600
$ws.Range("L126").Value = 18999.375
$ws.Range("M126").Value = -500002130
$ws.Range("N126").Value = -23939.375

# ---- Sheet index 7 ----
$ws = $wb.Worksheets.Item(7)
# Row 22
$ws.Range("H22").Value = 3485.8572
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 4680.2
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 4680.2
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -5270.2
# Row 27
$ws.Range("H27").Value = 3485.8572
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 4680.2
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 4680.2
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -4894.2
# Row 40
$ws.Range("H40").Value = 5705.5
$ws.Range("I40").Value = 2493.4
$ws.Range("K40").Value = 2493.4
$ws.Range("M40").Value = -2357.4
# Row 45
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
# Row 68
$ws.Range("H68").Value = 7494
$ws.Range("I68").Value = 7494
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 7494
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -6745
$ws.Range("N68").ClearContents()
# Row 71
$ws.Range("H71").Value = 7494
$ws.Range("I71").Value = 7494
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 37470
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -33726
$ws.Range("N71").ClearContents()
# Row 122
$ws.Range("H122").Value = 4402.6772
$ws.Range("I122").Value = 3606.1304
$ws.Range("K122").Value = 10818.3912
$ws.Range("M122").Value = -8368.3912
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()
# Row 136
$ws.Range("H136").Value = 14162.757
$ws.Range("I136").Value = 2886.1538
$ws.Range("K136").Value = 8658.4614
$ws.Range("M136").Value = -6108.4614
# Row 140
$ws.Range("H140").Value = 75585.60000000001
$ws.Range("J140").Value = 75585.60000000001
$ws.Range("L140").Value = 75585.60000000001
$ws.Range("N140").Value = -85945.60000000001

# ---- Sheet index 8 ----
$ws = $wb.Worksheets.Item(8)
# Row 52
$ws.Range("H52").Value = 2998.5
$ws.Range("I52").Value = 2998.5
$ws.Range("K52").Value = 2998.5
$ws.Range("M52").Value = -2772.5
# Row 54
$ws.Range("H54").Value = 27999.834
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
# Row 136
$ws.Range("H136").Value = 27301580
$ws.Range("I136").Value = 125001640
$ws.Range("J136").Value = 349839.62
$ws.Range("K136").Value = 375004920
$ws.Range("L136").Value = 1049518.86
$ws.Range("M136").Value = -375002370
$ws.Range("N136").Value = -1054618.86
